$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo: the date recorded as 2023-12-11 should have been 2023-12-31
$ws.Range("A64").Value = "2023-12-31 00:00:00"

# Append two new data points for 2024
$ws.Range("A65").Value = "2024-03-31 00:00:00"
$ws.Range("A66").Value = "2024-05-16 00:00:00"

# Copy A64's formatting (number format / style) onto the two new rows
$ws.Range("A64").Copy() | Out-Null
$ws.Range("A65:A66").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Move the active selection to the new last cell, matching the saved view state
$ws.Range("A65").Select() | Out-Null
